# ObservationDefinition profiles and examples updated for R5
#
# This script applies a set of shape-geometry nudges and text updates to
# slide 1 of the MedicationCatalog deck (ClinicalUseIssue -> ClinicalUse/
# ClinicalDefinition rename for R5).
#
# EMU <-> Point helper. PowerPoint's Shape.Left/Top/Width/Height and the
# xfrm off/ext the host eventually serializes are in points, while the
# diff we are replaying is expressed in EMU (1 pt = 12700 EMU). A tiny
# epsilon is added because the host's internal point<->EMU round trip
# truncates instead of rounding, which otherwise drops the low-order EMU
# about half the time.
function EMU([double]$emu) {
    return ($emu / 914400.0 * 72.0) + 0.00003
}

# Sets a paragraph's whole text while forcing the host's text-diff engine
# to emit a single run (it otherwise keeps the old run for any common
# prefix/suffix with the new text, which is still correct XML but not a
# literal match for the target). Flipping through an unrelated string
# first avoids any overlap with the old text.
function Set-ParaText($para, [string]$newText) {
    $para.Text = "~~~TMP~~~"
    $para.Text = $newText
}

# Same idea, but for a sub-range of characters located by locating the
# literal $oldText inside the shape's full text (1-based COM indexing).
function Set-SubText($textRange, [string]$oldText, [string]$newText) {
    $full = $textRange.Text
    $startIdx = $full.IndexOf($oldText) + 1
    $len = $oldText.Length
    $sub = $textRange.Characters($startIdx, $len)
    $sub.Text = "~~~TMP~~~"
    $sub2 = $textRange.Characters($startIdx, "~~~TMP~~~".Length)
    $sub2.Text = $newText
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function GetShape([int]$id) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# --- Rectangle 227 (id 228): nudge right ---
$sh = GetShape 228
$sh.Left = EMU 4766616

# --- Rectangle 228 (id 229): nudge right ---
$sh = GetShape 229
$sh.Left = EMU 3908856

# --- Rectangle 45 (id 46): move + widen, fix label text ---
$sh = GetShape 46
$sh.Left = EMU 3722533
$sh.Width = EMU 1333937
Set-ParaText $sh.TextFrame.TextRange.Paragraphs(2) "profile of ClinicalUseDefinition"

# --- Rectangle 141 (id 142): fix label text ---
$sh = GetShape 142
Set-ParaText $sh.TextFrame.TextRange.Paragraphs(1) "Indication profile of ClinicalUseDefinition"

# --- Rectangle 143 (id 144): move + widen, fix label text ---
$sh = GetShape 144
$sh.Left = EMU 2978150
$sh.Width = EMU 1176593
Set-ParaText $sh.TextFrame.TextRange.Paragraphs(2) "profile of ClinicalDefinition"

# --- Connecteur : en angle 153 (id 154): un-rotate/un-flip, reposition ---
$sh = GetShape 154
$sh.HorizontalFlip = 0
$sh.Rotation = 90
$sh.Left = EMU 2779682
$sh.Top = EMU 5186629
$sh.Width = EMU 1577748
$sh.Height = EMU 4218

# --- ZoneTexte 168 (id 169): fix label text (keep trailing "(4 slices)" run) ---
$sh = GetShape 169
Set-SubText $sh.TextFrame.TextRange "ClinicalUseIssue " "ClinicalUseDefinition "

# --- ZoneTexte 177 (id 178): nudge right ---
$sh = GetShape 178
$sh.Left = EMU 4372003

# --- Connecteur : en angle 233 (id 234): reposition/resize ---
$sh = GetShape 234
$sh.Left = EMU 3123028
$sh.Top = EMU 6750879
$sh.Width = EMU 263876
$sh.Height = EMU 1404492

# --- Connecteur : en angle 236 (id 237): reposition/resize ---
$sh = GetShape 237
$sh.Left = EMU 3188319
$sh.Top = EMU 6703797
$sh.Width = EMU 574075
$sh.Height = EMU 1828293

# --- Connecteur : en angle 240 (id 241): reposition/resize ---
$sh = GetShape 241
$sh.Left = EMU 3227380
$sh.Top = EMU 6646527
$sh.Width = EMU 912932
$sh.Height = EMU 2262252

# --- Connecteur : en angle 132 (id 133): reposition/resize + adjust handle ---
$sh = GetShape 133
$sh.Left = EMU 2815294
$sh.Top = EMU 5155234
$sh.Width = EMU 2329578
$sh.Height = EMU 818837
$sh.Adjustments.Item(1) = 33373 / 100000.0

# --- Rectangle 151 (id 152): widen, fix label text ---
$sh = GetShape 152
$sh.Width = EMU 1295504
Set-ParaText $sh.TextFrame.TextRange.Paragraphs(2) "profile of ClinicalUseDefinition"

# --- Connecteur : en angle 152 (id 153): reposition/resize + adjust handle ---
$sh = GetShape 153
$sh.Left = EMU 2261438
$sh.Top = EMU 5709090
$sh.Width = EMU 4263571
$sh.Height = EMU 1645117
$sh.Adjustments.Item(1) = 18277 / 100000.0

# --- ZoneTexte 156 (id 157): nudge right ---
$sh = GetShape 157
$sh.Left = EMU 4099701

# --- ZoneTexte 164 (id 165): reposition (both x & y) ---
$sh = GetShape 165
$sh.Left = EMU 1983473
$sh.Top = EMU 8934554

# --- Rectangle 167 (id 168): move + widen, fix label text ---
$sh = GetShape 168
$sh.Left = EMU 5370088
$sh.Width = EMU 1295504
Set-ParaText $sh.TextFrame.TextRange.Paragraphs(2) "profile of ClinicalUseDefinition"

# --- Connecteur : en angle 169 (id 170): reposition/resize + adjust handle ---
$sh = GetShape 170
$sh.Left = EMU 2240436
$sh.Top = EMU 5730092
$sh.Width = EMU 5107632
$sh.Height = EMU 2447175
$sh.Adjustments.Item(1) = 15189 / 100000.0

# --- Connecteur : en angle 181 (id 182): shrink width ---
$sh = GetShape 182
$sh.Width = EMU 4028144
